$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "partner_type"
$ws.Range("J1").Value = "self_journal_id"

$ws.Range("I4").Value = "other"
$ws.Range("J4").Value = "z0bug.sale"

$ws.Range("J4").Select()
